$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Campaign Report date" value (I2): 2018-07-31 to 2018-08-23 -> 2018-07-31 to 2018-08-26
$ws.Range("I2").Value = "2018-07-31 to 2018-08-26"

# Update row 10 VDX delivery numbers
$ws.Range("L10").Value = 126934
$ws.Range("M10").Value = 0.1353962666666667
$ws.Range("N10").Value = 20309.4423673191
